$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 6666.6665
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5586

$ws.Range("H133").Value = 27998
$ws.Range("J133").Value = 27998
$ws.Range("L133").Value = 27998
$ws.Range("N133").Value = -38118

$ws.Range("H134").Value = 26992.5
$ws.Range("J134").Value = 26992.5
$ws.Range("L134").Value = 26992.5
$ws.Range("N134").Value = -37132.5

$ws.Range("H136").Value = 57390
$ws.Range("J136").Value = 57390
$ws.Range("L136").Value = 57390
$ws.Range("N136").Value = -67590

$ws.Range("H137").Value = 2332199
$ws.Range("I137").Value = 6067151.5
$ws.Range("K137").Value = 18201454.5
$ws.Range("M137").Value = -18198904.5

$ws.Range("H139").Value = 40497.5
$ws.Range("J139").Value = 43996.668
$ws.Range("L139").Value = 43996.668
$ws.Range("N139").Value = -54276.668

$ws.Range("H140").Value = 48266.668
$ws.Range("J140").Value = 48266.668
$ws.Range("L140").Value = 48266.668
$ws.Range("N140").Value = -58626.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 50000224
$ws.Range("I8").Value = 100000000
$ws.Range("J8").Value = 450
$ws.Range("K8").Value = 100000000
$ws.Range("L8").Value = 450
$ws.Range("M8").Value = -99999856
$ws.Range("N8").Value = -738

$ws.Range("H32").Value = 1370188.5
$ws.Range("I32").Value = 1416606.8
$ws.Range("J32").Value = 850
$ws.Range("K32").Value = 1416606.8
$ws.Range("L32").Value = 850
$ws.Range("M32").Value = -1416319.8
$ws.Range("N32").Value = -1424

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 720
$ws.Range("I5").Value = 366.66666
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 366.66666
$ws.Range("L5").Value = 1250
$ws.Range("M5").Value = -253.66666
$ws.Range("N5").Value = -1476

$ws.Range("H10").Value = 102.5
$ws.Range("I10").Value = 102.5
$ws.Range("K10").Value = 102.5
$ws.Range("M10").Value = 37.5

$ws.Range("H24").Value = 1027.5
$ws.Range("I24").Value = 433
$ws.Range("K24").Value = 433
$ws.Range("M24").Value = -198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 17500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 17500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 17500
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -17840

$ws.Range("H19").Value = 41934.168
$ws.Range("I19").Value = 321
$ws.Range("J19").Value = 250000
$ws.Range("K19").Value = 321
$ws.Range("L19").Value = 250000
$ws.Range("M19").Value = -151
$ws.Range("N19").Value = -250340

$ws.Range("H24").Value = 41934.168
$ws.Range("I24").Value = 321
$ws.Range("J24").Value = 250000
$ws.Range("K24").Value = 321
$ws.Range("L24").Value = 250000
$ws.Range("M24").Value = -151
$ws.Range("N24").Value = -250340

$ws.Range("H31").Value = 1012470.56
$ws.Range("I31").Value = 859.0238000000001
$ws.Range("J31").Value = 1897630.8
$ws.Range("K31").Value = 859.0238000000001
$ws.Range("L31").Value = 1897630.8
$ws.Range("M31").Value = -564.0238000000001
$ws.Range("N31").Value = -1898220.8

$ws.Range("H34").Value = 1012470.56
$ws.Range("I34").Value = 859.0238000000001
$ws.Range("J34").Value = 1897630.8
$ws.Range("K34").Value = 859.0238000000001
$ws.Range("L34").Value = 1897630.8
$ws.Range("M34").Value = -657.0238000000001
$ws.Range("N34").Value = -1898034.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.56
$ws.Range("I2").Value = 130.9
$ws.Range("J2").Value = 45.333332
$ws.Range("K2").Value = 785.4000000000001
$ws.Range("L2").Value = 271.999992
$ws.Range("M2").Value = -672.4000000000001
$ws.Range("N2").Value = -497.999992

$ws.Range("H7").Value = 362.58823
$ws.Range("I7").Value = 134.83333
$ws.Range("J7").Value = 486.81818
$ws.Range("K7").Value = 404.49999
$ws.Range("L7").Value = 1460.45454
$ws.Range("M7").Value = -292.49999
$ws.Range("N7").Value = -1684.45454

$ws.Range("H17").Value = 10001
$ws.Range("I17").Value = 10001
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 30003
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -29834
$ws.Range("N17").ClearContents()

$ws.Range("H23").Value = 436.85715
$ws.Range("I23").Value = 57
$ws.Range("J23").Value = 466.07693
$ws.Range("K23").Value = 171
$ws.Range("L23").Value = 1398.23079
$ws.Range("M23").Value = 64
$ws.Range("N23").Value = -1868.23079

$ws.Range("H34").Value = 952.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 952.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2857.5
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3025.5

$ws.Range("H39").Value = 5300
$ws.Range("J39").Value = 5300
$ws.Range("L39").Value = 15900
$ws.Range("N39").Value = -16488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2781246.2
$ws.Range("I132").Value = 4313547.5
$ws.Range("K132").Value = 12940642.5
$ws.Range("M132").Value = -12938112.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4250
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3388
$ws.Range("N3").Value = -5224

$ws.Range("H15").Value = 4250
$ws.Range("I15").Value = 3500
$ws.Range("J15").Value = 5000
$ws.Range("K15").Value = 3500
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = -3330
$ws.Range("N15").Value = -5340

$ws.Range("H136").Value = 3394.8367
$ws.Range("I136").Value = 1728.0588
$ws.Range("J136").Value = 7172.8667
$ws.Range("K136").Value = 5184.1764
$ws.Range("L136").Value = 21518.6001
$ws.Range("M136").Value = -2634.1764
$ws.Range("N136").Value = -26618.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 518.6667
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H133").Value = 31475.715
$ws.Range("J133").Value = 31475.715
$ws.Range("L133").Value = 31475.715
$ws.Range("N133").Value = -41595.715
